# Renal cancer.xlsx -- "Refined metadata to be additional tab"
#
# 1. Update the panel_query_time timestamps on the "data" sheet (column F).
# 2. Add a new "metadata" worksheet (after "data") describing the panel
#    query itself (data_name / data_id / data_version / ... columns).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Refresh the per-row query timestamps on the data sheet.
# ---------------------------------------------------------------------
$ws1.Range("F2").Value  = "2021-10-05 14:35:33.260904"
$ws1.Range("F3").Value  = "2021-10-05 14:35:33.260912"
$ws1.Range("F4").Value  = "2021-10-05 14:35:33.260916"
$ws1.Range("F5").Value  = "2021-10-05 14:35:33.260919"
$ws1.Range("F6").Value  = "2021-10-05 14:35:33.260921"
$ws1.Range("F7").Value  = "2021-10-05 14:35:33.260924"
$ws1.Range("F8").Value  = "2021-10-05 14:35:33.260927"
$ws1.Range("F9").Value  = "2021-10-05 14:35:33.260929"
$ws1.Range("F10").Value = "2021-10-05 14:35:33.260932"
$ws1.Range("F11").Value = "2021-10-05 14:35:33.260935"

# ---------------------------------------------------------------------
# 2. Add the "metadata" worksheet right after "data".
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Re-use the existing header style (bold font + border + centered,
# style index 1 in the source workbook) by copying formats from cells
# that already carry it, instead of building a brand new style.
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("B1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Renal cancer"
$ws2.Range("C2").Value = 3278
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.3"
$ws2.Range("D2").ClearFormats()
$ws2.Range("E2").Value = "2020-08-10T07:05:57.075400Z"
$ws2.Range("F2").Value = "2021-10-05 14:35:33.257030"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3278/?format=json"

# Keep "data" as the active sheet/selection, matching the source file
# (the commit only appends a sheet entry; bookViews are untouched).
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null

Write-Output "metadata sheet added; timestamps refreshed"
